# "add test file with vexing characters in worksheet names"
#
# Inserts a new worksheet "Gabe's S''heet" before the existing "HiGabe!!!!"
# sheet, gives it the same data ("numbers" header + 1..5), adds sheet-scoped
# (local) defined names "A_one"/"numbers" pointing at the new sheet (in
# addition to the existing workbook-scoped ones that still point at
# "HiGabe!!!!"), and leaves "HiGabe!!!!" as the active sheet with cell E11
# selected.

$wb = $excel.ActiveWorkbook
$hi = $wb.Worksheets.Item("HiGabe!!!!")

# New sheet goes before "HiGabe!!!!" (matches sheetId="2"/rId1 ordering).
$new = $wb.Worksheets.Add($hi)
$new.Name = "Gabe's S''heet"

# Same contents as "HiGabe!!!!": a "numbers" header followed by 1..5.
# The "numbers" text reuses the existing shared string.
$new.Range("A1").Value = "numbers"
$new.Range("A2").Value = 1
$new.Range("A3").Value = 2
$new.Range("A4").Value = 3
$new.Range("A5").Value = 4
$new.Range("A6").Value = 5

# Sheet-scoped defined names on the new sheet (localSheetId="0"), mirroring
# the workbook-scoped "A_one"/"numbers" that already refer to "HiGabe!!!!".
# Sheet names containing a literal apostrophe need it doubled once the name
# is wrapped in the surrounding quotes of a sheet-qualified reference.
$escapedName = $new.Name.Replace("'", "''")
$refOne = "='" + $escapedName + "'!`$A`$1"
$refNumbers = "='" + $escapedName + "'!`$A`$2:`$A`$6"
$new.Names.Add("A_one", $refOne) | Out-Null
$new.Names.Add("numbers", $refNumbers) | Out-Null

# Re-fetch "HiGabe!!!!" (the reference captured before Add() goes stale)
# and make it the active sheet again with E11 selected.
$hi = $wb.Worksheets.Item("HiGabe!!!!")
$hi.Activate()
$hi.Range("E11").Select() | Out-Null
